# mansoni_coverage_scenario_1.xlsx - "Updates to prop never treated and coverage times"
#
# 1) Platform Coverage sheet, row 2: the "proportion never treated"/coverage
#    values used to be entered only every other year (H,J,L,N,...). Fill in
#    the skipped odd columns (I,K,M,O,Q,S,U,W,Y,AA,AC) with the same 0.6
#    value so coverage is specified for every year H2:AD2.
# 2) Update the saved view/selection on the "Platform Coverage" sheet and
#    the workbook window to reflect where the author left the cursor
#    (scrolled right to show the new columns, selection on AB5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Platform Coverage")

# --- 1) Fill in the previously-empty columns on row 2 -----------------
foreach ($col in @("I","K","M","O","Q","S","U","W","Y","AA","AC")) {
    $ws.Range($col + "2").Value = 0.6
}

# --- 2) Update sheet view / selection ----------------------------------
$ws.Activate()

$win = $excel.ActiveWindow
$win.ScrollColumn = 22   # topLeftCell -> V1
$win.ScrollRow = 1

$ws.Range("AB5").Select()

# --- Update workbook window size/position -------------------------------
$appWin = $wb.Windows.Item(1)
$appWin.Left = -110
$appWin.Top = -110
$appWin.Width = 19420
$appWin.Height = 10300
